# Added range of motion display
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Writing")

# New values for the existing rows (column F added)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 10

# New rows showing the range-of-motion settings
$ws.Range("A4").Value = "Max_RoM"
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = "Min_RoM"
$ws.Range("B5").Value = -20

$ws.Range("A6").Value = "Period"
$ws.Range("B6").Value = 2

# Update the selection / view to match the new state
$ws.Range("J10").Select()
